# Converts the DELAWARE_2023 MCAS sheet:
#  - renames header columns to snake_case machine-readable names
#  - converts the ALL-CAPS state/municipality text in columns A and B to Proper Case
#  - removes the trailing footnote/metadata rows (489-493)

function ToProperCase($s) {
    $lower = $s.ToLower()
    $upper = $s.ToUpper()
    $result = ""
    $prevIsLetter = $false
    for ($i = 0; $i -lt $s.Length; $i++) {
        $ch = $lower.Substring($i, 1)
        if ($ch -match "\p{L}") {
            if ($prevIsLetter) {
                $result = $result + $ch
            } else {
                $result = $result + $upper.Substring($i, 1)
            }
            $prevIsLetter = $true
        } else {
            $result = $result + $ch
            $prevIsLetter = $false
        }
    }
    return $result
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row to machine-friendly snake_case names
$ws.Cells.Item(1,1).Value2 = "mx_state"
$ws.Cells.Item(1,2).Value2 = "mx_municipality"
$ws.Cells.Item(1,3).Value2 = "n_matriculas"
$ws.Cells.Item(1,4).Value2 = "pct_matriculas"

# 2. Proper-case the state (A) and municipality (B) text, rows 2..487
$lastDataRow = 487
$range = $ws.Range("A2:B" + $lastDataRow)
$vals = $range.Value2

for ($r = 1; $r -le $vals.GetLength(0); $r++) {
    for ($c = 1; $c -le $vals.GetLength(1); $c++) {
        $v = $vals[$r, $c]
        if ($v -ne $null) {
            $vals[$r, $c] = ToProperCase($v)
        }
    }
}

$range.Value2 = $vals

# 3. Remove the trailing footnote / metadata rows (488-493)
$ws.Range("A488:A493").EntireRow.Delete()
